# Daily Update 키워드 10개
# Remove the literal "<b>" / "</b>" highlighting markup that was left in the
# product-title column (B) of the worksheet, e.g.
#   "바퀴달린 이동식 간이 쇼파 <b>테이블</b> 미니 ..."
# becomes
#   "바퀴달린 이동식 간이 쇼파 테이블 미니 ..."
# Only column B (rows 2-101) contains this markup; all other columns and
# the header row are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $value = $cell.Value2

    if ($value -ne $null -and $value -match '<b>|</b>') {
        $newValue = $value -replace '<b>', '' -replace '</b>', ''
        $cell.Value2 = $newValue
    }
}
